# Fruta / hortaliza, semanal
# Insert a new weekly price-log row at row 105 (pushing existing rows
# 105..157 down to 106..158) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, shifting rows 105..157 down to 106..158.
# -4121 = xlShiftDown
$ws.Rows.Item(105).Insert(-4121)

# Populate the newly inserted row 105 with this week's observation.
$ws.Cells.Item(105, 1).Value2  = 4
$ws.Cells.Item(105, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(105, 3).Value2  = "Los Lagos"
$ws.Cells.Item(105, 4).Value2  = 44510
$ws.Cells.Item(105, 5).Value2  = 10
$ws.Cells.Item(105, 6).Value2  = "Fruta"
$ws.Cells.Item(105, 7).Value2  = 100108
$ws.Cells.Item(105, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(105, 9).Value2  = 100108005
$ws.Cells.Item(105, 10).Value2 = "Piña"
$ws.Cells.Item(105, 11).Value2 = "Caramelo"
$ws.Cells.Item(105, 12).Value2 = "Segunda"
$ws.Cells.Item(105, 13).Value2 = 200
$ws.Cells.Item(105, 14).Value2 = 22000
$ws.Cells.Item(105, 15).Value2 = 23000
$ws.Cells.Item(105, 16).Value2 = 22500
$ws.Cells.Item(105, 17).Value2 = "`$/caja 14 unidades"
$ws.Cells.Item(105, 18).Value2 = "Ecuador"
$ws.Cells.Item(105, 19).Value2 = 1607
$ws.Cells.Item(105, 20).Value2 = 14
